# Auto-generated Excel COM-interop edit script
# Applies numeric 'want to go' count updates and a cancellation notice
# to the performance (天鹅湖) event, mirroring the target commit diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览 (Exhibition)
$ws2 = $wb.Worksheets.Item(2)   # 演出 (Performance)
$ws3 = $wb.Worksheets.Item(3)   # 本地生活 (Local Life)
$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (All Types)


# --- Sheet 1 ---
$ws1.Range("F3").Value = 154
$ws1.Range("F4").Value = 1787
$ws1.Range("F6").Value = 1106
$ws1.Range("F7").Value = 2235
$ws1.Range("F8").Value = 2147
$ws1.Range("F9").Value = 1119
$ws1.Range("F11").Value = 23
$ws1.Range("F13").Value = 404
$ws1.Range("F15").Value = 44
$ws1.Range("F16").Value = 305
$ws1.Range("F17").Value = 234
$ws1.Range("F18").Value = 1602
$ws1.Range("F19").Value = 8
$ws1.Range("F20").Value = 651
$ws1.Range("F21").Value = 739
$ws1.Range("F22").Value = 11
$ws1.Range("F23").Value = 621
$ws1.Range("F24").Value = 12319
$ws1.Range("F25").Value = 12373
$ws1.Range("F26").Value = 917
$ws1.Range("F27").Value = 708
$ws1.Range("F29").Value = 247
$ws1.Range("F30").Value = 29
$ws1.Range("F31").Value = 389
$ws1.Range("F34").Value = 7
$ws1.Range("F35").Value = 209
$ws1.Range("F36").Value = 606

# --- Sheet 2 ---
$ws2.Range("C3").Value = "广州·柴可夫斯基百年经典《天鹅湖》室内乐重奏音乐会（取消）"
$ws2.Range("G3").Value = "不可售"
$ws2.Range("F7").Value = 45

# --- Sheet 3 ---
$ws3.Range("F2").Value = 78
$ws3.Range("F3").Value = 54

# --- Sheet 4 ---
$ws4.Range("F3").Value = 78
$ws4.Range("F4").Value = 154
$ws4.Range("F5").Value = 1787
$ws4.Range("F7").Value = 1106
$ws4.Range("F8").Value = 2235
$ws4.Range("F9").Value = 2147
$ws4.Range("F10").Value = 1119
$ws4.Range("F12").Value = 54
$ws4.Range("F13").Value = 23
$ws4.Range("F15").Value = 404
$ws4.Range("F18").Value = 44
$ws4.Range("C19").Value = "广州·柴可夫斯基百年经典《天鹅湖》室内乐重奏音乐会（取消）"
$ws4.Range("G19").Value = "不可售"
$ws4.Range("F20").Value = 305
$ws4.Range("F22").Value = 234
$ws4.Range("F23").Value = 1602
$ws4.Range("F24").Value = 9
$ws4.Range("F25").Value = 651
$ws4.Range("F26").Value = 739
$ws4.Range("F27").Value = 11
$ws4.Range("F28").Value = 621
$ws4.Range("F29").Value = 12319
$ws4.Range("F30").Value = 12373
$ws4.Range("F31").Value = 917
$ws4.Range("F32").Value = 708
$ws4.Range("F34").Value = 247
$ws4.Range("F35").Value = 29
$ws4.Range("F36").Value = 389
$ws4.Range("F40").Value = 7
$ws4.Range("F42").Value = 209
$ws4.Range("F43").Value = 606
$ws4.Range("F44").Value = 45
